# Mark GroupBy (row 14), All (row 50), and Contains (row 52) as "Done"
# by setting column A on their rows, matching the existing "Done" marker
# used throughout the sheet for completed operators.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A14").Value = "Done"
$ws.Range("A50").Value = "Done"
$ws.Range("A52").Value = "Done"

# Update the view state to match where the user ended up scrolled/selected.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("D20").Select()
